# Price parser custom activity - update Parts/Category naming + replace Webshop list with hyperlinks

$wb = $excel.ActiveWorkbook

# --- Parts sheet: shorten a few component names (categories stay the same) ---
$wsParts = $wb.Worksheets.Item("Parts")
$wsParts.Range("A3").Value = "B650 AORUS ELITE AX AMD B650"
$wsParts.Range("A5").Value = "ATX Tower Gaming"
$wsParts.Range("A6").Value = "GeForce RTX 4060"

# --- Category sheet: capitalize "Hard drive" ---
$wsCategory = $wb.Worksheets.Item("Category")
$wsCategory.Range("A3").Value = "Hard drive "

# --- Webshop sheet: replace shop list with new shops + real hyperlinks ---
$wsWebshop = $wb.Worksheets.Item("Webshop")

$wsWebshop.Range("A2").Value = "Emag"
$wsWebshop.Range("B2").Value = "http://emag.hu"

$wsWebshop.Range("A3").Value = "Alza"
$wsWebshop.Range("B3").Value = "http://alza.hu"

$wsWebshop.Range("A4").Value = "Pcx"
$wsWebshop.Range("B4").Value = "http://pcx.hu"

$wsWebshop.Hyperlinks.Add($wsWebshop.Range("B2"), "http://emag.hu")
$wsWebshop.Hyperlinks.Add($wsWebshop.Range("B3"), "http://alza.hu")
$wsWebshop.Hyperlinks.Add($wsWebshop.Range("B4"), "http://pcx.hu")

# --- Restore selection / active cells on each sheet ---
$wsCategory.Range("A4").Select() | Out-Null
$wsWebshop.Range("A5").Select() | Out-Null
$wsParts.Activate() | Out-Null
$wsParts.Range("F2").Select() | Out-Null
